$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New style (cellXfs index 43): fontId=2, fillId=0, borderId=0, horizontal=center ---
# Built by copying an existing "plain/no-border" format (M22) and centering it.
# Used by M32.
$ws.Range("M22").Copy()
$ws.Range("M32").PasteSpecial(-4122)          # xlPasteFormats
$ws.Range("M32").HorizontalAlignment = -4108  # xlCenter

# --- New style (cellXfs index 44): fontId=1 (bold), fillId=2, borderId=0, horizontal=center ---
# Built by copying an existing bold/filled header format (B1), then stripping its border
# and resetting vertical alignment back to the (unwritten) default.
# Used by M31.
$ws.Range("B1").Copy()
$ws.Range("M31").PasteSpecial(-4122)          # xlPasteFormats
$ws.Range("M31").Borders.LineStyle = -4142    # xlLineStyleNone
$ws.Range("M31").VerticalAlignment = -4107    # xlVAlignBottom (default -> not serialized)

# --- Row 31: K31, L31, M31 ---
$ws.Range("K20").Copy()
$ws.Range("K31").PasteSpecial(-4122)
$ws.Range("K31").Value = "Test Case Name"

$ws.Range("K20").Copy()
$ws.Range("L31").PasteSpecial(-4122)
$ws.Range("L31").Value = "INC Code"

# --- Row 32: K32, L32, M32 ---
# New shared string order matches the source diff: "Setup Configuration" (140)
# is introduced before "Attribute Group" (141).
$ws.Range("K21").Copy()
$ws.Range("K32").PasteSpecial(-4122)
$ws.Range("K32").Value = "Setup Configuration"

$ws.Range("M31").Value = "Attribute Group"

$ws.Range("L29").Copy()
$ws.Range("L32").PasteSpecial(-4122)
$ws.Range("L32").Value = "00016"

$ws.Range("M32").Value = "AttrGroupTest123"

# --- Row 32 height ---
$ws.Rows.Item(32).RowHeight = 12.75

# --- Selection / view ---
$ws.Range("L32").Select()
